function Set-CellText($sheet, $addr, $text) {
    $c = $sheet.Range($addr)
    $c.Value = "'" + $text
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "335.78"
Set-CellText $ws "E2" "1.77%"
Set-CellText $ws "D3" "44.10"
Set-CellText $ws "E3" "6.55%"
Set-CellText $ws "D4" "5.754"
Set-CellText $ws "E4" "1.83%"
Set-CellText $ws "D5" "0.08390"
Set-CellText $ws "E5" "2.09%"
Set-CellText $ws "D6" "8.857"
Set-CellText $ws "E6" "1.14%"
Set-CellText $ws "D7" "1.966"
Set-CellText $ws "E7" "-2.10%"
Set-CellText $ws "E8" "-3.09%"
Set-CellText $ws "D9" "0.9503"
Set-CellText $ws "E9" "2.67%"
Set-CellText $ws "D10" "0.1247"
Set-CellText $ws "E10" "-2.17%"
Set-CellText $ws "D11" "0.1976"
Set-CellText $ws "E11" "1.03%"
Set-CellText $ws "D12" "0.1032"
Set-CellText $ws "E12" "10.14%"
Set-CellText $ws "D13" "0.04435"
Set-CellText $ws "E13" "13.39%"
Set-CellText $ws "E14" "0.61%"
Set-CellText $ws "D15" "0.001287"
Set-CellText $ws "E15" "-1.20%"
Set-CellText $ws "D16" "0.006014"
Set-CellText $ws "E16" "-1.86%"
Set-CellText $ws "D17" "3.494"
Set-CellText $ws "E17" "1.39%"
Set-CellText $ws "D18" "4.520"
Set-CellText $ws "E18" "-0.36%"
Set-CellText $ws "D19" "0.3537"
Set-CellText $ws "E19" "1.66%"
Set-CellText $ws "D20" "8.673"
Set-CellText $ws "E20" "4.34%"
Set-CellText $ws "E21" "-0.78%"
Set-CellText $ws "E22" "-0.77%"
Set-CellText $ws "E23" "0.44%"
Set-CellText $ws "D24" "0.001258"
Set-CellText $ws "E24" "0.11%"
Set-CellText $ws "D25" "0.004364"
Set-CellText $ws "E25" "1.06%"
Set-CellText $ws "E26" "5.20%"
Set-CellText $ws "D27" "0.0003996"
Set-CellText $ws "E27" "-94.67%"
Set-CellText $ws "D39" "0.02831"
Set-CellText $ws "E39" "2.40%"
Set-CellText $ws "D40" "0.05988"
Set-CellText $ws "E40" "8.39%"
Set-CellText $ws "D41" "0.007933"
Set-CellText $ws "E41" "0.18%"
Set-CellText $ws "D42" "0.1429"
Set-CellText $ws "E42" "0.49%"
Set-CellText $ws "D43" "0.008965"
Set-CellText $ws "E43" "0.22%"
Set-CellText $ws "E44" "0.29%"
Set-CellText $ws "D45" "0.01017"
Set-CellText $ws "E45" "-14.30%"
Set-CellText $ws "D46" "0.00007274"
Set-CellText $ws "E46" "3.85%"
Set-CellText $ws "E47" "0.19%"
Set-CellText $ws "D48" "0.003203"
Set-CellText $ws "E48" "0.35%"
Set-CellText $ws "D49" "0.002274"
Set-CellText $ws "E49" "-0.29%"
Set-CellText $ws "E50" "0.19%"
Set-CellText $ws "E51" "0.19%"

Write-Host "Applied 68 cell updates"
